$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '27.956.50'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +3.26%  '
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '1.724.99'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +2.92%  '
$ws.Range("E4").Value = '  -0.26%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '218.78'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +1.58%  '
$ws.Range("E6").Value = '  +1.37%  '
$ws.Range("E7").Value = '  -0.24%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '23.98'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +12.41%  '
$ws.Range("E9").Value = '  +3.64%  '
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '0.0634'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  +1.87%  '
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.0902'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  +2.13%  '
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '1.969.10'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  +2.95%  '
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '1.728.77'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +3.25%  '
$ws.Range("E14").Value = '  +3.55%  '
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '0.567'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +5.46%  '
$cell = $ws.Range("D16")
$cell.NumberFormat = "@"
$cell.Value = '67.96'
$cell.Style = "Normal"
$ws.Range("E16").Value = '  +2.77%  '
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '27.902.87'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  +3.14%  '
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '244.20'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  +2.80%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0756'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  +2.32%  '
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '7.90'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -3.08%  '
$ws.Range("E21").Value = '  -0.17%  '
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = '4.65'
$cell.Style = "Normal"
$ws.Range("E22").Value = '  +3.93%  '
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '9.74'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +4.74%  '
$ws.Range("E24").Value = '  +0.90%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '149.36'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +1.30%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '7.52'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  +3.83%  '
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '16.82'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  +2.58%  '
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '0.115'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +1.83%  '
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("E30").Value = '  +2.85%  '
$ws.Range("E31").Value = '  +1.85%  '
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '3.45'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +2.91%  '
$ws.Range("E33").Value = '  +3.21%  '
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '1.488.68'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -2.80%  '
$ws.Range("E35").Value = '  -1.93%  '
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.614'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  +3.40%  '
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '0.960'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  +4.54%  '
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = '2.40'
$cell.Style = "Normal"
$ws.Range("E38").Value = '  +0.40%  '
$ws.Range("E39").Value = '  +0.43%  '
$ws.Range("E40").Value = '  -0.91%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '71.45'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +5.50%  '
$ws.Range("E42").Value = '  +5.93%  '
$ws.Range("E43").Value = '  -0.23%  '
$ws.Range("B44").Value = 'MXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '2.29'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +1.11%  '
$ws.Range("B45").Value = 'RocketPoolETH'
$ws.Range("C45").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '1.873.11'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  +2.92%  '
$ws.Range("E46").Value = '  +1.36%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '1.74'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +11.99%  '
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '91.17'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  +0.49%  '
$ws.Range("E49").Value = '  +2.93%  '
$ws.Range("E50").Value = '  +3.08%  '
$ws.Range("E51").Value = '  +1.27%  '
